$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the brand new shared-string-backed cells first, in the specific
# order that reproduces the target shared string table ordering.
$ws.Range("C24").Value = "value_attribute_classify"
$ws.Range("D24").Value = "Table_Field / cell-df (with VA classification)"
$ws.Range("D28").Value = "cell_analysis"
$ws.Range("C28").Value = "analyze_cells / analyse_cells"
$ws.Range("C25").Value = "basic_classifier"
$ws.Range("C26").Value = "numeric_values_classifier"
$ws.Range("C27").Value = "sample_based_classifier"
$ws.Range("E24").Value = "general VA classifier (which takes VA classification function as input)"
$ws.Range("E26").Value = "number like cells are considered as values and rest as atribute. "
$ws.Range("E25").Value = "basic VA classifier which consider numeric (already numeric) as values and text as attribute"
$ws.Range("E27").Value = "sample based VA classifier where connected cells are marked same as the match."
$ws.Range("E28").Value = "complete heuristic based structure allocation"
$ws.Range("C29").Value = "compose_cells"
$ws.Range("D29").Value = "tibble"
$ws.Range("E29").Value = "tidy form (without datablock wise column collation)"
$ws.Range("C30").Value = "compose_cells(discard_raw_cols = TRUE)"
$ws.Range("E30").Value = "tidy form (without datablock wise column collation) (ideal if data-block wise variation is not present)"

# Fill the cells that reuse already-existing shared strings.
$ws.Range("B25").Value = "Table_Field / cell-df"
$ws.Range("D25").Value = "Table_Field / cell-df (with VA classification)"
$ws.Range("B26").Value = "Table_Field / cell-df"
$ws.Range("D26").Value = "Table_Field / cell-df (with VA classification)"
$ws.Range("B27").Value = "Table_Field / cell-df"
$ws.Range("D27").Value = "Table_Field / cell-df (with VA classification)"
$ws.Range("B28").Value = "Table_Field / cell-df (with VA classification)"
$ws.Range("B29").Value = "cell_analysis"
$ws.Range("B30").Value = "cell_analysis"
$ws.Range("D30").Value = "tibble"

# Styles: rows 24-28 (all columns) use the same style as the existing body
# rows (style applied to row 23). Columns C:E on rows 29-30 keep the
# bottom-border style that used to sit on row 26 (captured first, before
# row 26 gets overwritten by the body-style paste below); column B on rows
# 29-30 uses the plain body style instead.
$ws.Range("C26:E26").Copy() | Out-Null
$ws.Range("C29:E30").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B23:E23").Copy() | Out-Null
$ws.Range("B24:E28").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("B23").Copy() | Out-Null
$ws.Range("B29:B30").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Column widths (targets: 46.28515625 / 44 / 106.7109375 "characters"; the
# engine quantizes ColumnWidth to the nearest 1/6 of a character when it
# round-trips through pixels, so feed it the input that lands on the
# closest achievable grid point).
$ws.Columns.Item(2).ColumnWidth = 45.5
$ws.Columns.Item(3).ColumnWidth = 43.166666666666664
$ws.Columns.Item(5).ColumnWidth = 105.83333333333333

# View: scroll so column B is the left-most visible column, and move the
# active selection to D15.
$ws.Range("D15").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
